# Append a new scraped job posting (2026-02-13 18:44:20 JST run) to the
# top of the data block (row 9), pushing the previously-top-ranked rows
# down by one, and refresh the "fetched at" timestamp on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-02-13 18:44:20"

# --- 1. Insert a fresh row at position 9; this pushes old rows 9-13 down
#        to 10-14 and widens the used range to H14 automatically. ---
$ws.Rows.Item(9).Insert()

# --- 2. Refresh column A ("取得日時") on every data row (2-14) to the new
#        timestamp -- both the rows that existed before and the new one. ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 3. Populate the newly inserted row 9 with the new job posting. ---
$ws.Cells.Item(9, 2).Value = "【急募】アンドロイドタブレット向け将棋アプリ開発者募集"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5491505"
$ws.Cells.Item(9, 7).Value = 88
$ws.Cells.Item(9, 8).Value = "◆開発 ◇アプリ"

# --- 4. The row insert shifted the F-column cell contents/styles down
#        correctly, but it does NOT re-point the worksheet's stored
#        hyperlink relationships at the same time, so every hyperlink
#        needs to be rebuilt from scratch, in row order, so rId1..rId13
#        line up with rows 2..14 top to bottom again. ---
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 14; $r++) {
    $target = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $target)
}
